$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 05:16"

# Update country rows with refreshed case data (values & reordered countries)
$rows = @(
  @{Row=22; Values=@("Australia", 1895, 8, 118, 1769, 11, 1, 8)}
  @{Row=33; Values=@("Tailandia", 827, 106, 52, 771, 7, 3, 4)}
  @{Row=34; Values=@("Polonia", 749, 0, 13, 728, 3, 0, 8)}
  @{Row=35; Values=@("Chile", 746, 0, 11, 733, 7, 0, 2)}
  @{Row=46; Values=@("India", 500, 1, 37, 453, 0, 0, 10)}
  @{Row=54; Values=@("Hong Kong", 357, 0, 100, 253, 4, 0, 4)}
  @{Row=107; Values=@("Trinidad yTobago", 52, 1, 0, 52, 0, 0, 0)}
  @{Row=108; Values=@("Liechtenstein", 51, 0, 0, 51, 0, 0, 0)}
  @{Row=110; Values=@("Afganistan", 42, 2, 1, 40, 0, 0, 1)}
  @{Row=111; Values=@("Cuba", 40, 0, 0, 39, 3, 0, 1)}
  @{Row=114; Values=@("Mauricio", 36, 0, 0, 34, 1, 0, 2)}
  @{Row=115; Values=@("Consejo Danes para los Refugiados", 36, 0, 0, 34, 0, 0, 2)}
  @{Row=123; Values=@("Paraguay", 27, 5, 0, 25, 1, 1, 2)}
  @{Row=124; Values=@("Costa de Marfil", 25, 0, 2, 23, 0, 0, 0)}
  @{Row=125; Values=@("Macao", 25, 0, 10, 15, 0, 0, 0)}
  @{Row=126; Values=@("Mayotte", 24, 0, 0, 24, 0, 0, 0)}
  @{Row=127; Values=@("Monaco", 23, 0, 1, 22, 0, 0, 0)}
  @{Row=128; Values=@("Guyana", 20, 0, 0, 19, 0, 0, 1)}
  @{Row=129; Values=@("Guatemala", 20, 0, 0, 19, 0, 0, 1)}
  @{Row=143; Values=@("Aruba", 12, 3, 1, 11, 0, 0, 0)}
  @{Row=144; Values=@("Etiopia", 11, 0, 0, 11, 0, 0, 0)}
  @{Row=145; Values=@("Mongolia", 10, 0, 0, 10, 0, 0, 0)}
  @{Row=146; Values=@("Guinea Ecuatorial", 9, 0, 0, 9, 0, 0, 0)}
  @{Row=147; Values=@("Uganda", 9, 0, 0, 9, 0, 0, 0)}
  @{Row=151; Values=@("Surinam", 6, 1, 0, 6, 0, 0, 0)}
  @{Row=154; Values=@("Bermudas", 6, 0, 0, 6, 0, 0, 0)}
  @{Row=155; Values=@("Gabon", 6, 0, 0, 5, 0, 0, 1)}
  @{Row=156; Values=@("El Salvador", 5, 2, 0, 5, 0, 0, 0)}
)

foreach ($r in $rows) {
  $rowNum = $r.Row
  $vals = $r.Values
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item($rowNum, $col).Value = $vals[$i]
  }
}

Write-Output "Done updating paises worksheet"
